$d = $word.ActiveDocument

# Insert a new paragraph after the last one ("Elizabeth"), matching the
# same run/paragraph formatting (lang="es-ES"), and give it the text
# "Brandon".
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Brandon"

# The "_GoBack" bookmark (added by Word whenever a document is edited and
# saved) now belongs at the new last-edited location, i.e. around the
# newly typed "Brandon" text, instead of around "Elizabeth". Re-adding a
# bookmark with the same reserved name moves it off the old paragraph and
# onto the new one.
$d.Bookmarks.Add("_GoBack", $newPara.Range)
